# "Uno mas para subir" - append a couple more lines under "OtroDoc" and
# mark the non-dictionary words ("OtroDoc" and "dfasdf") the same way
# Word's background spell-checker would, with <w:proofErr> spellStart/
# spellEnd pairs around the offending run.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Flag "OtroDoc" (the existing paragraph) as a spelling error ---
# Rewrite just the run content of paragraph 1, leaving the paragraph
# mark itself (and its paraId/textId/rsid attributes) untouched.
$p1 = $d.Paragraphs(1)
$p1Range = $p1.Range
$p1Text = $d.Range($p1Range.Start, $p1Range.End - 1)
[void]$p1Text.InsertXML("<w:p $wNs><w:proofErr w:type=`"spellStart`"/><w:r><w:t>OtroDoc</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>")

# --- 2. Append the new paragraphs at the end of the document body ---
# One blank line, then "s", then "dfasdf" (also flagged as misspelled).
$bodyEnd = $d.Content.End
$tail = $d.Range($bodyEnd - 1, $bodyEnd - 1)
$newParas = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document $wNs>
        <w:body>
          <w:p/>
          <w:p><w:r><w:t>s</w:t></w:r></w:p>
          <w:p><w:proofErr w:type="spellStart"/><w:r><w:t>dfasdf</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
[void]$tail.InsertXML($newParas)
